$d = $word.ActiveDocument

# Helper: locate the paragraph whose text contains $needle.
function Get-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- Change 1: "Review Project proposal draft" -------------------------
# Split the single run into two runs around "draft" and wrap "draft" with
# proofErr gramStart/gramEnd markers (mirrors Word's grammar-checker
# splitting the run at a flagged word boundary).
$p1 = Get-ParagraphByText $d "Review Project proposal draft"
if ($p1 -ne $null) {
    $full1 = $p1.Range
    $xml1 = '<w:p w14:paraId="7835BDF3" w14:textId="6E213DDF" w:rsidR="00392303" w:rsidRPr="009E53A3" w:rsidRDefault="00392303">' +
            '<w:r w:rsidRPr="009E53A3"><w:t xml:space="preserve">Review Project proposal </w:t></w:r>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r w:rsidRPr="009E53A3"><w:t>draft</w:t></w:r>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '</w:p>'
    $full1.InsertXML($xml1)
}

# --- Change 2: "... before 30-10-2023" ----------------------------------
# Split the trailing " 30-10-2023" run into a run containing just the
# space and a second run containing the date, wrapping the date with
# proofErr gramStart/gramEnd markers.
$p2 = Get-ParagraphByText $d "30-10-2023"
if ($p2 -ne $null) {
    $full2 = $p2.Range
    $xml2 = '<w:p w14:paraId="424BA1C8" w14:textId="4EB58FA4" w:rsidR="00B532BE" w:rsidRPr="00BE56BB" w:rsidRDefault="00114509" w:rsidP="00BE56BB">' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
            '<w:r w:rsidRPr="00BE56BB"><w:t xml:space="preserve">Draft of the proposal written </w:t></w:r>' +
            '<w:r w:rsidR="005C0D4F" w:rsidRPr="00BE56BB"><w:t xml:space="preserve">to be </w:t></w:r>' +
            '<w:r w:rsidR="00390114" w:rsidRPr="00BE56BB"><w:t>sent</w:t></w:r>' +
            '<w:r w:rsidR="005C0D4F" w:rsidRPr="00BE56BB"><w:t xml:space="preserve"> before</w:t></w:r>' +
            '<w:r w:rsidR="00643F08" w:rsidRPr="00BE56BB"><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r w:rsidR="00643F08" w:rsidRPr="00BE56BB"><w:t>30-10-2023</w:t></w:r>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '</w:p>'
    $full2.InsertXML($xml2)
}

Write-Output "Done."
